$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3875
$ws.Range("B2").Value = 2034
$ws.Range("C2").Value = 4578
$ws.Range("D2").Value = 3876
$ws.Range("E2").Value = 3399
$ws.Range("F2").Value = 4636
